# "auto advance in player" -------------------------------------------------
# The Player (row 8) and Player_S16 (row 9) states both now react to the
# Play/Up/Down buttons the same way, and the old Playing/Paused sub-states
# are replaced by a spacer row + a single "auto advance" marker row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$btnUp   = "Button (LV_EVENT_SHORT_CLICKED)`nEVT_BUTTON_UP_CLICKED"
$btnDown = "Button (LV_EVENT_SHORT_CLICKED)`nEVT_BUTTON_DOWN_CLICKED"
$btnSet  = "Button (LV_EVENT_SHORT_CLICKED)`nEVT_BUTON_SETTING_CLICKED"
$btnBack = "Button (LV_EVENT_SHORT_CLICKED)`nEVT_BUTTON_BACK_CLICKED"
$btnPlay = "Button (LV_EVENT_SHORT_CLICKED)`nEVT_BUTTON_PLAY_CLICKED"

# --- Row 8 ("Player"): add the new Up / Down click columns -----------------
$ws.Range("E8").WrapText = $true
$ws.Range("E8").VerticalAlignment = -4160   # xlVAlignTop
$ws.Range("E8").Value = $btnUp

$ws.Range("F8").WrapText = $true
$ws.Range("F8").VerticalAlignment = -4160   # xlVAlignTop
$ws.Range("F8").Value = $btnDown

# --- Row 9 ("Player_S16"): gains the same Setting/Back/Play/Up/Down cells
#     the Player row has -------------------------------------------------
$ws.Range("B9").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("B9").WrapText = $true
$ws.Range("B9").HorizontalAlignment = -4131 # xlVAlignLeft
$ws.Range("B9").Value = $btnSet

$ws.Range("C9").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("C9").WrapText = $true
$ws.Range("C9").Value = $btnBack

$ws.Range("D9").WrapText = $true
$ws.Range("D9").Value = $btnPlay

$ws.Range("E9").WrapText = $true
$ws.Range("E9").VerticalAlignment = -4160   # xlVAlignTop
$ws.Range("E9").Value = $btnUp

$ws.Range("F9").WrapText = $true
$ws.Range("F9").VerticalAlignment = -4160   # xlVAlignTop
$ws.Range("F9").Value = $btnDown

# --- Row 10: the old "Player_S16_Playing" sub-state row becomes a blank
#     spacer row --------------------------------------------------------
$ws.Range("A10").Value = ""

# --- Row 11: the old "Player_S16_Paused" sub-state row is replaced by a
#     single auto-advance marker -----------------------------------------
$ws.Range("A11").Value = "∂"

# --- row heights: every multi-line (wrapped, 2-row-tall) row now renders
#     slightly taller than before ----------------------------------------
$ws.Rows.Item(3).RowHeight = 32
$ws.Rows.Item(4).RowHeight = 32
$ws.Rows.Item(8).RowHeight = 32
$ws.Rows.Item(9).RowHeight = 32
$ws.Rows.Item(14).RowHeight = 32
$ws.Rows.Item(16).RowHeight = 32

# --- cosmetic view state to match the refreshed sheet --------------------
$ws.Columns.Item(1).ColumnWidth = 19.6
$ws.Columns.Item(7).ColumnWidth = 26

$ws.Select()
$ws.Range("A11").Select()
$excel.ActiveWindow.Zoom = 130
